# PROS-10763 - MARSRU - New KPIs
#
# Adds 4 new KPI rows (ДМП - display stand related questions) to both the
# "KPI with Codes" and "KPI with Names" sheets, mirroring the existing
# "INSERT INTO kpi_level_2 ... VALUES(...)" row layout.

$wb = $excel.ActiveWorkbook

# Literal fragments re-used by every generated SQL row. They all contain a
# leading apostrophe, and Excel's plain `.Value = "'..."` assignment treats a
# leading apostrophe as a "force text" prefix and silently drops it - so
# these are written through a formula (`="literal"`) and then flattened back
# to a plain value/shared-string via copy + paste-values, which sidesteps
# that stripping.
$litF = "('"
$litG = "', NULL, '"
$litH = "', '"
$litI = "', '20', '3', '5', 'Custom', '3', '1', '0', '0', '0', '0', '0',"
$litJ = "),"

function Set-LiteralText($ws, $addr, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
}

function Add-KpiRow-Codes($ws, $row, $a, $b, $c, $e) {
    # Column A keeps the "s=27" number style used by the rest of this block;
    # copy it in from the last pre-existing row of the same family (151).
    $ws.Range("A151").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $a
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Formula = "=IF(EXACT(C$row,""Boolean""),""'10'"",""NULL"")"
    $ws.Range("E$row").Value = $e
    Set-LiteralText $ws "F$row" $litF
    Set-LiteralText $ws "G$row" $litG
    Set-LiteralText $ws "H$row" $litH
    Set-LiteralText $ws "I$row" $litI
    Set-LiteralText $ws "J$row" $litJ
    $ws.Range("K$row").Formula = "=CONCATENATE(F$row,E$row,G$row,A$row,H$row,A$row,I$row,D$row,J$row)"
}

function Add-KpiRow-Names($ws, $row, $a, $b, $c, $e) {
    $ws.Range("A151").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $a
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Formula = "=IF(EXACT(C$row,""Boolean""),""'10'"",""NULL"")"
    $ws.Range("E$row").Value = $e
    Set-LiteralText $ws "F$row" $litF
    Set-LiteralText $ws "G$row" $litG
    Set-LiteralText $ws "H$row" $litH
    Set-LiteralText $ws "I$row" $litI
    Set-LiteralText $ws "J$row" $litJ
    $ws.Range("K$row").Formula = "=CONCATENATE(F$row,E$row,G$row,CONCATENATE(A$row,""-RUS""),H$row,CONCATENATE(A$row,"" - "",B$row),I$row,D$row,J$row)"
}

$typeDmp = "Выберите тип ДМП"
$promo = "Проходит ли с данного дисплея промо активность?"
$facings = "Укажите общее количество фейсингов Марс на ДМП (кол-во фейсингов из цены мотивационная программа)"

# ---- Sheet 1: "KPI with Codes" ----
$ws1 = $wb.Worksheets.Item("KPI with Codes")

Add-KpiRow-Codes $ws1 153 10011 $typeDmp  "Text"    1437
Add-KpiRow-Codes $ws1 154 10012 $promo    "Boolean" 1438
Add-KpiRow-Codes $ws1 155 10014 $facings  "Int"     1439
Add-KpiRow-Codes $ws1 156 10013 $typeDmp  "Text"    1440

# Flatten the F:J helper-formula cells for the new rows into literal shared
# strings (matches how the rest of the sheet stores these constants).
$ws1.Range("F153:J156").Copy()
$ws1.Range("F153:J156").PasteSpecial(-4163)

$ws1.Columns.Item(2).ColumnWidth = 19

$ws1.Activate() | Out-Null
$ws1.Range("D156").Select() | Out-Null

# ---- Sheet 2: "KPI with Names" ----
$ws2 = $wb.Worksheets.Item("KPI with Names")

Add-KpiRow-Names $ws2 153 10011 $typeDmp  "Text"    2732
Add-KpiRow-Names $ws2 154 10012 $promo    "Boolean" 2733
Add-KpiRow-Names $ws2 155 10014 $facings  "Int"     2734
Add-KpiRow-Names $ws2 156 10013 $typeDmp  "Text"    2735

$ws2.Range("F153:J156").Copy()
$ws2.Range("F153:J156").PasteSpecial(-4163)

$ws2.Columns.Item(2).ColumnWidth = 19

$ws2.Activate() | Out-Null
$ws2.Range("B156").Select() | Out-Null

$ws1.Activate() | Out-Null
